$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching style of existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-26
$values = @{
    2  = @(8, 8)
    3  = @(7, 8)
    4  = @(8, 8)
    5  = @(7, 7)
    6  = @(1, 2)
    7  = @(7, 7)
    8  = @(1, 1)
    9  = @(7, 7)
    10 = @(5, 6)
    11 = @(1, 1)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(1, 2)
    15 = @(1, 2)
    16 = @(8, 9)
    17 = @(8, 8)
    18 = @(6, 6)
    19 = @(6, 7)
    20 = @(9, 9)
    21 = @(7, 8)
    22 = @(6, 7)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(2, 2)
    26 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
